$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 12 with data mirroring row 6's client (456), but a new
# cargo_id (4565) and a new current-status text with trailing period.
# Write the new shared-string values in the same order the target
# workbook introduces them: status text first, then the cargo id.
$ws.Range("H12").Value = "Поступил на склад в Китай."
$ws.Range("C12").Value = "4565"

$ws.Range("B12").Value = 456
$ws.Range("D12").Value = 22
$ws.Range("E12").Value = 22
$ws.Range("F12").Value = "2"
$ws.Range("G12").Value = 2
$ws.Range("I12").Value = "02,02,2020"
$ws.Range("J12").Value = "Отправлен в РФ"
$ws.Range("K12").Value = "03,02,2020"
$ws.Range("L12").Value = "Таможенное оформление"
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = "Прибыл в РФ"
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = "Выдан получателю"
$ws.Range("Q12").Value = 0

# Match the number formats of the analogous cells in row 6 (client/cargo id
# columns carry the numeric/text custom formats applied there).
$ws.Range("B12").NumberFormat = $ws.Range("B6").NumberFormat
$ws.Range("C12").NumberFormat = $ws.Range("C6").NumberFormat
$ws.Range("F12").NumberFormat = $ws.Range("F6").NumberFormat
$ws.Range("I12").NumberFormat = $ws.Range("I6").NumberFormat

# Move the active selection as recorded in the edited workbook.
$ws.Range("D14").Select()
